# Update "想去人数" (wanted-to-go count) figures in column F on the
# "展览" (Exhibition) and "全部类型" (All types) sheets, reflecting the
# latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (first sheet) ---
$wsExhibition = $wb.Worksheets.Item(1)
$wsExhibition.Range("F4").Value = 750
$wsExhibition.Range("F11").Value = 8267
$wsExhibition.Range("F12").Value = 441
$wsExhibition.Range("F20").Value = 10384
$wsExhibition.Range("F35").Value = 857
$wsExhibition.Range("F39").Value = 2880
$wsExhibition.Range("F41").Value = 152
$wsExhibition.Range("F44").Value = 32
$wsExhibition.Range("F46").Value = 95

# --- Sheet "全部类型" (fourth sheet) ---
$wsAllTypes = $wb.Worksheets.Item(4)
$wsAllTypes.Range("F4").Value = 750
$wsAllTypes.Range("F11").Value = 8267
$wsAllTypes.Range("F12").Value = 441
$wsAllTypes.Range("F20").Value = 10384
$wsAllTypes.Range("F34").Value = 857
$wsAllTypes.Range("F38").Value = 2880
$wsAllTypes.Range("F41").Value = 152
$wsAllTypes.Range("F44").Value = 32
$wsAllTypes.Range("F46").Value = 95
